$d = $word.ActiveDocument

# 1) Update the "Status da anomalia" field text in the first table
$d.Content.Find.Execute("Aprovada para resolução", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Reparada e encerrada.", 2)

# 2) Fill in the last (empty) row of the revision-history table (2nd table)
$tbl = $d.Tables.Item(2)
$lastRow = $tbl.Rows.Count

$values = @("10/06/2015", "3", "Alteração do status e encerramento", "Moisés")
for ($col = 1; $col -le 4; $col++) {
    $cellRange = $tbl.Cell($lastRow, $col).Range
    $cellRange.Text = $values[$col - 1]
    $cellRange.Font.Name = "Times New Roman"
    $cellRange.Font.Bold = $true
}
